# Update cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.807.19"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "'1.872.32"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'0.7315"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'241.27"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3137"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.07126"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'24.40"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.08151"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").Value = "'1.908.98"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "'0.7418"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'92.37"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'29.830.80"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'6.005"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "'248.37"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'13.39"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'0.000007798"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "'2.159.98"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'7.756"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").Value = "'0.1541"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "'9.205"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "'164.12"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "'2.018"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'1.448"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'4.518"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'1.522"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'4.181"
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("D34").Value = "'0.05310"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'1.230"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "'0.7407"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'2.703"
$ws.Range("D39").Value = "'0.01935"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'2.735"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'0.4464"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").Value = "'5.966"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").Value = "'0.8677"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "'71.25"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "'1.045.62"
$ws.Range("E45").Value = "  -5.92%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'103.96"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'1.817"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").Value = "'7.431"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.054.63"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.512"
$ws.Range("E51").Value = "  -0.39%  "
